# The document's auto-maintained "_GoBack" bookmark (Word drops one at the
# location of the most recent edit) currently sits between the two runs
# " On average, e" / "ach session takes about 45 minutes." in the Procedure
# paragraph of Experiment 1. The edit moves the user's last-edit location to
# a new (empty) paragraph right after the "Results" heading of Experiment 1,
# so the bookmark must be removed from its old spot and re-created there.

$d = $word.ActiveDocument

# 1) Drop the old _GoBack bookmark from the "On average, e|ach session" spot.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2) Locate the "Results" heading that introduces the Experiment 1 results
#    section (the first "Results" heading in the paper).
$found = $d.Content
[void]$found.Find.Execute("Results", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$resultsPara = $found.Paragraphs(1)

# 3) Press Enter at the end of that heading to open a new, empty paragraph
#    right after it.
$resultsPara.Range.InsertParagraphAfter()
$headingIndex = $resultsPara.Range.Paragraphs(1).Index
$newPara = $d.Paragraphs($headingIndex + 1)

# 4) Give the new (still empty) paragraph its final shape: no explicit
#    paragraph style (falls back to the document default), an East Asian
#    font hint stamped on the paragraph mark, and the relocated _GoBack
#    bookmark collapsed at that same empty paragraph.
$newParaXml = '<?xml version="1.0" standalone="yes"?>' + `
  '<?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p>' + `
              '<w:pPr>' + `
                '<w:rPr>' + `
                  '<w:rFonts w:eastAsia="MS Mincho" w:hint="eastAsia"/>' + `
                '</w:rPr>' + `
              '</w:pPr>' + `
              '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
              '<w:bookmarkEnd w:id="0"/>' + `
            '</w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

[void]$newPara.Range.InsertXML($newParaXml)
